# Updated line loss / power-flow results for the 380 kV case (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","E","F","G","H","I","K","M")

$data = @(
    @(0.264701291173111, 0.04279265066249138, 0.155992084214823, 0.4443680307746263, 0.7875958996127252, 0.8637624282434473, 0.7105404212837456, 0.2842951474065387, 0.2289909675444406),
    @(0.2346120983195021, 0.03727357622452132, 0.1449832165851177, 0.387822817061874, 0.7885799687407768, 0.8690837821736892, 0.7157964830893597, 0.2503307836205408, 0.2071061116803037),
    @(0.2161603653072177, 0.03387878703541958, 0.1383270589604066, 0.3531389305168915, 0.7897599676313831, 0.8727805951642864, 0.719439602606009, 0.2294884120931044, 0.1937621308272242),
    @(0.2086472397703574, 0.03249383540565987, 0.135640358407116, 0.3390132514313251, 0.7903851682388137, 0.8743949677764959, 0.7210285844773061, 0.2209981818728579, 0.188347693198601),
    @(0.2074000681861889, 0.0322637710716549, 0.1351957811587923, 0.336668177824194, 0.7904976872494558, 0.8746695464281657, 0.7212987323782727, 0.2195885856972666, 0.1874500356190651),
    @(0.2160590155694138, 0.03386011539078027, 0.1382907212895361, 0.3529483938344953, 0.7897678154615448, 0.8728019304524679, 0.7194606097968759, 0.2293738966098857, 0.1936890155029616),
    @(0.2543218917618049, 0.0408909138260185, 0.1521746284958923, 0.4248636149813336, 0.7878154362507672, 0.8655080286595762, 0.7122662789933258, 0.2725818865574468, 0.2214255871762987),
    @(0.3295304497594884, 0.05463205085544587, 0.1802334015684792, 0.5661985755041457, 0.7885762470052242, 0.8546175250931611, 0.7014670059342905, 0.3574018098875626, 0.27656759910505),
    @(0.3848872835358748, 0.06470296864520719, 0.2013752515792078, 0.6702781546542269, 0.791963333727665, 0.8487043598333059, 0.6955626887565103, 0.4197730153111365, 0.3175555041579301),
    @(0.4100916983894081, 0.069279945352406, 0.2111120281015744, 0.7176906081379002, 0.7941249120798517, 0.8464693945928587, 0.6933202219017218, 0.4481594019385682, 0.3363091507948823),
    @(0.4196389717394879, 0.07101254905212784, 0.2148165455445863, 0.7356546913071611, 0.7950332227183736, 0.8456886280309703, 0.6925350483562482, 0.4589104068303413, 0.3434264401367386),
    @(0.4175826709899866, 0.07063942803608825, 0.2140179328242482, 0.7317853510981394, 0.7948336013416224, 0.8458538617148434, 0.6927012993756989, 0.4565949137827943, 0.3418929042273362),
    @(0.4108771017314155, 0.06942249954128954, 0.2114164505087572, 0.7191683204515869, 0.7941978373915504, 0.8464038455900749, 0.6932543412032928, 0.4490438609021226, 0.3368943791253187),
    @(0.4067701168568476, 0.06867701857677844, 0.2098252428379439, 0.7114413442032514, 0.7938201189086698, 0.8467492695175878, 0.6936014369384935, 0.4444188370607094, 0.3338346861220529),
    @(0.3832405390132578, 0.06440376737785414, 0.2007413523871904, 0.6671810134426437, 0.7918346039193267, 0.8488595873071887, 0.6957181823424321, 0.4179181481342766, 0.3163320964599166),
    @(0.3688114054590699, 0.06178117381394088, 0.1951994085708293, 0.6400460337125793, 0.7907759067857256, 0.8502708275154447, 0.6971304784891146, 0.4016641026146033, 0.3056226059129372),
    @(0.3605142674568071, 0.06027231785286347, 0.192023036887278, 0.6244449056556647, 0.7902253725651036, 0.8511253502392009, 0.6979845181965345, 0.3923164816415294, 0.2994729504139428),
    @(0.3577053720627248, 0.05976137405809823, 0.1909494889271812, 0.619163680173358, 0.7900489861968225, 0.8514220260719725, 0.6982808404279197, 0.3891517641051507, 0.2973925230933432),
    @(0.3703471943414343, 0.06206039557746124, 0.1957881960739272, 0.642933953830422, 0.7908825590515249, 0.8501161663886023, 0.6969758171434179, 0.4033942431316859, 0.306761596814745),
    @(0.4128466134091298, 0.06977995696605888, 0.2121800942682981, 0.7228739723491628, 0.794382136407549, 0.8462405213995936, 0.6930901606020683, 0.4512617439689564, 0.3383621400454189),
    @(0.4406392563889199, 0.07482166283425329, 0.2229947585348953, 0.7751780083420101, 0.7971927807730594, 0.844089800727275, 0.6909238093098438, 0.4825557052014631, 0.359106450130291),
    @(0.4258043594211074, 0.07213111888734147, 0.2172133834976009, 0.7472568307830727, 0.7956446240055328, 0.8452026586014654, 0.692045811017401, 0.4658527098900436, 0.348026395207512),
    @(0.3696528692655647, 0.06193416282299324, 0.1955219748487949, 0.6416283278902171, 0.7908341605400153, 0.8501859541787127, 0.697045608437044, 0.4026120557897457, 0.3062466360777378),
    @(0.3091663897014598, 0.05091930051403892, 0.172551619601812, 0.5279251897347166, 0.7878759686038137, 0.857197552630538, 0.7040329468727506, 0.3344462768597793, 0.2615680635631819)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellRef = "$($cols[$j])$row"
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}
